$d = $word.ActiveDocument

function Get-ParaIndexByText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }
    $rngFresh = $d.Range($rng.Start, $rng.End)
    $rngFresh.Expand(4) | Out-Null   # wdParagraph -> expand to the whole paragraph
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Start -eq $rngFresh.Start -and $pp.Range.End -eq $rngFresh.End) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "Field Types: Textual Data" paragraph -> expanded field-type notes
# ---------------------------------------------------------------------------
$idx1 = Get-ParaIndexByText("Field Types: Textual Data")
if ($idx1 -gt 0) {
    $r1 = $d.Paragraphs.Item($idx1).Range
    $xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t xml:space="preserve">Field Types: Textual Data </w:t></w:r><w:r><w:rPr/><w:t>(email, text, characters)</w:t></w:r><w:r><w:rPr/><w:br/><w:t xml:space="preserve">Field Types: </w:t></w:r><w:r><w:rPr/><w:t>Numeric</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> Data </w:t></w:r><w:r><w:rPr/><w:t>(integers, decimals)</w:t></w:r><w:r><w:rPr/><w:br/><w:t xml:space="preserve">Field Types: </w:t></w:r><w:r><w:rPr/><w:t>MIscellaneous</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> Data </w:t></w:r><w:r><w:rPr/><w:t>(boolean, dateTIme)</w:t><w:br/><w:t>Field Types: Relational Data (foreignkey, manytomanyfield)</w:t><w:br/><w:t xml:space="preserve">FIELDS contain attributes) </w:t><w:br/><w:tab/><w:t>blank(attribute)=true ==not required)</w:t><w:br/><w:tab/><w:t xml:space="preserve">null(attribute) (no data) | </w:t><w:br/><w:tab/><w:t>choices(atribute) (set of choices)</w:t></w:r></w:p>'
    $r1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# 2) "CREAR EL MODELO DE LOS DATOS QUE LLEVARA" paragraph -> bold
# ---------------------------------------------------------------------------
$idx2 = Get-ParaIndexByText("CREAR EL MODELO")
if ($idx2 -gt 0) {
    $r2 = $d.Paragraphs.Item($idx2).Range
    $xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:b/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="es-MX"/></w:rPr><w:t>CREAR EL MODELO DE LOS DATOS QUE LLEVAR&#193;</w:t></w:r></w:p>'
    $r2.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# 4) New paragraph "after MAKINGS MIGRATIONS you need to MIGRATE" right
#    after the "there are default apps..." paragraph. (Done before (3) so
#    the lower-index insertion in (3) cannot perturb this resolved index.)
# ---------------------------------------------------------------------------
$idx4 = Get-ParaIndexByText("there are default apps when making the initial migration, they com with models and migrations")
if ($idx4 -gt 0) {
    $p4 = $d.Paragraphs.Item($idx4)
    $p4.Range.InsertParagraphAfter()
    $p4b = $d.Paragraphs.Item($idx4 + 1)
    $p4b.Range.Text = "after MAKINGS MIGRATIONS you need to MIGRATE"
}

# ---------------------------------------------------------------------------
# 3) New paragraph "(to make the initial migration)" right before the
#    "python manage.py makemigrations" paragraph.
# ---------------------------------------------------------------------------
$idx3 = Get-ParaIndexByText("migration created for a new Django app will create tables for the models that are defined")
if ($idx3 -gt 0) {
    $p3 = $d.Paragraphs.Item($idx3)
    $p3.Range.InsertParagraphAfter()
    $p3b = $d.Paragraphs.Item($idx3 + 1)
    $p3b.Range.Text = "(to make the initial migration)"
}

Write-Output "done"
